# Auto-generated edit script: update FFXIV leve-profit market price columns (H-N)
# across all 8 job sheets, per scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 459.5
$ws.Range("I4").Value = 320.5
$ws.Range("K4").Value = 320.5
$ws.Range("M4").Value = -206.5
$ws.Range("H40").Value = 3349.8
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H138").Value = 2362.9666
$ws.Range("I138").Value = 2544.1875
$ws.Range("J138").Value = 2155.8572
$ws.Range("K138").Value = 7632.5625
$ws.Range("L138").Value = 6467.571599999999
$ws.Range("M138").Value = -2492.5625
$ws.Range("N138").Value = -16747.5716

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3020.6
$ws.Range("I2").Value = 916.8461
$ws.Range("K2").Value = 916.8461
$ws.Range("M2").Value = -803.8461
$ws.Range("H61").Value = 16763.666
$ws.Range("I61").Value = 3529.4285
$ws.Range("J61").Value = 25185.455
$ws.Range("K61").Value = 3529.4285
$ws.Range("L61").Value = 25185.455
$ws.Range("M61").Value = -3317.4285
$ws.Range("N61").Value = -25609.455
$ws.Range("H74").Value = 26644.059
$ws.Range("I74").Value = 3737.2856
$ws.Range("K74").Value = 3737.2856
$ws.Range("M74").Value = -2863.2856
$ws.Range("H77").Value = 26644.059
$ws.Range("I77").Value = 3737.2856
$ws.Range("K77").Value = 18686.428
$ws.Range("M77").Value = -14318.428
$ws.Range("H97").Value = 6096.933
$ws.Range("I97").Value = 2788.1667
$ws.Range("K97").Value = 2788.1667
$ws.Range("M97").Value = -2292.1667
$ws.Range("H116").Value = 3020.6
$ws.Range("I116").Value = 916.8461
$ws.Range("K116").Value = 916.8461
$ws.Range("M116").Value = 1377.1539
$ws.Range("H132").Value = 4013196.2
$ws.Range("I132").Value = 4443.95
$ws.Range("K132").Value = 13331.85
$ws.Range("M132").Value = -10801.85
$ws.Range("H136").Value = 16763.666
$ws.Range("I136").Value = 3529.4285
$ws.Range("J136").Value = 25185.455
$ws.Range("K136").Value = 10588.2855
$ws.Range("L136").Value = 75556.36500000001
$ws.Range("M136").Value = -8038.2855
$ws.Range("N136").Value = -80656.36500000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4977.647
$ws.Range("I94").Value = 5235.826
$ws.Range("J94").Value = 4437.8184
$ws.Range("K94").Value = 5235.826
$ws.Range("L94").Value = 4437.8184
$ws.Range("M94").Value = -4784.826
$ws.Range("N94").Value = -5339.8184
$ws.Range("H99").Value = 1549.8823
$ws.Range("I99").Value = 1175.5454
$ws.Range("J99").Value = 2236.1667
$ws.Range("K99").Value = 1175.5454
$ws.Range("L99").Value = 2236.1667
$ws.Range("M99").Value = 322.4546
$ws.Range("N99").Value = -5232.1667
$ws.Range("H105").Value = 2153.7693
$ws.Range("I105").Value = 1501.125
$ws.Range("K105").Value = 1501.125
$ws.Range("M105").Value = 245.875
$ws.Range("H106").Value = 17390.666
$ws.Range("J106").Value = 17390.666
$ws.Range("L106").Value = 17390.666
$ws.Range("N106").Value = -19914.666
$ws.Range("H107").Value = 2207
$ws.Range("I107").Value = 1985.1765
$ws.Range("J107").Value = 2835.5
$ws.Range("K107").Value = 1985.1765
$ws.Range("L107").Value = 2835.5
$ws.Range("M107").Value = -65.17650000000003
$ws.Range("N107").Value = -6675.5
$ws.Range("H134").Value = 11524.615
$ws.Range("I134").Value = 6232.72
$ws.Range("K134").Value = 18698.16
$ws.Range("M134").Value = -16163.16

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1205.75
$ws.Range("I22").Value = 936.0769
$ws.Range("K22").Value = 936.0769
$ws.Range("M22").Value = -586.0769
$ws.Range("H31").Value = 23106
$ws.Range("I31").Value = 8416.5
$ws.Range("K31").Value = 8416.5
$ws.Range("M31").Value = -8121.5
$ws.Range("H34").Value = 23106
$ws.Range("I34").Value = 8416.5
$ws.Range("K34").Value = 8416.5
$ws.Range("M34").Value = -8214.5
$ws.Range("H62").Value = 4330.6
$ws.Range("I62").Value = 3751
$ws.Range("J62").Value = 5200
$ws.Range("K62").Value = 3751
$ws.Range("L62").Value = 5200
$ws.Range("M62").Value = -3127
$ws.Range("N62").Value = -6448
$ws.Range("H65").Value = 4330.6
$ws.Range("I65").Value = 3751
$ws.Range("J65").Value = 5200
$ws.Range("K65").Value = 18755
$ws.Range("L65").Value = 26000
$ws.Range("M65").Value = -15635
$ws.Range("N65").Value = -32240
$ws.Range("H107").Value = 3412.08
$ws.Range("I107").Value = 1754.6428
$ws.Range("K107").Value = 1754.6428
$ws.Range("M107").Value = 165.3571999999999
$ws.Range("H122").Value = 4863.963
$ws.Range("I122").Value = 2032.1765
$ws.Range("J122").Value = 9678
$ws.Range("K122").Value = 6096.529500000001
$ws.Range("L122").Value = 29034
$ws.Range("M122").Value = -3646.529500000001
$ws.Range("N122").Value = -33934
$ws.Range("H132").Value = 8589.75
$ws.Range("I132").Value = 2681.8333
$ws.Range("K132").Value = 8045.499899999999
$ws.Range("M132").Value = -5515.499899999999
$ws.Range("H134").Value = 31256538
$ws.Range("I134").Value = 1833.579
$ws.Range("J134").Value = 76936500
$ws.Range("K134").Value = 5500.737
$ws.Range("L134").Value = 230809500
$ws.Range("M134").Value = -2965.737
$ws.Range("N134").Value = -230814570
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1000
$ws.Range("I8").Value = 1000
$ws.Range("K8").Value = 3000
$ws.Range("M8").Value = -2861
$ws.Range("H120").Value = 12500
$ws.Range("I120").Value = 12500
$ws.Range("K120").Value = 37500
$ws.Range("M120").Value = -32662
$ws.Range("H131").Value = 1488.02
$ws.Range("J131").Value = 1488.02
$ws.Range("L131").Value = 4464.059999999999
$ws.Range("N131").Value = -14544.06
$ws.Range("H139").Value = 7581.1963
$ws.Range("I139").Value = 11534.1875
$ws.Range("J139").Value = 6000
$ws.Range("K139").Value = 34602.5625
$ws.Range("L139").Value = 18000
$ws.Range("M139").Value = -29462.5625
$ws.Range("N139").Value = -28280

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 7746.893
$ws.Range("I126").Value = 5875.3125
$ws.Range("K126").Value = 17625.9375
$ws.Range("M126").Value = -15155.9375
$ws.Range("H132").Value = 19551.234
$ws.Range("I132").Value = 19551.234
$ws.Range("K132").Value = 58653.702
$ws.Range("M132").Value = -56123.702
$ws.Range("H135").Value = 154018.16
$ws.Range("J135").Value = 154018.16
$ws.Range("L135").Value = 154018.16
$ws.Range("N135").Value = -164158.16

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1491.6774
$ws.Range("I16").Value = 1279.8636
$ws.Range("K16").Value = 1279.8636
$ws.Range("M16").Value = -1109.8636
$ws.Range("H46").Value = 5725.8335
$ws.Range("J46").Value = 6999.875
$ws.Range("L46").Value = 6999.875
$ws.Range("N46").Value = -7375.875
$ws.Range("H55").Value = 2013.5
$ws.Range("I55").Value = 1034.3
$ws.Range("K55").Value = 1034.3
$ws.Range("M55").Value = -861.3
$ws.Range("H132").Value = 762260
$ws.Range("I132").Value = 2782.3667
$ws.Range("K132").Value = 8347.1001
$ws.Range("M132").Value = -5817.1001
$ws.Range("H136").Value = 13117
$ws.Range("I136").Value = 12659.3
$ws.Range("J136").Value = 13727.267
$ws.Range("K136").Value = 37977.89999999999
$ws.Range("L136").Value = 41181.801
$ws.Range("M136").Value = -35427.89999999999
$ws.Range("N136").Value = -46281.801
$ws.Range("H140").Value = 110654.5
$ws.Range("J140").Value = 110654.5
$ws.Range("L140").Value = 110654.5
$ws.Range("N140").Value = -121014.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2024.8334
$ws.Range("I96").Value = 1429.6
$ws.Range("K96").Value = 1429.6
$ws.Range("M96").Value = -56.59999999999991
$ws.Range("H113").Value = 2334.425
$ws.Range("J113").Value = 2360.0833
$ws.Range("L113").Value = 7080.249899999999
$ws.Range("N113").Value = -11420.2499
$ws.Range("H126").Value = 5455.304
$ws.Range("I126").Value = 4511.357
$ws.Range("J126").Value = 6923.6665
$ws.Range("K126").Value = 13534.071
$ws.Range("L126").Value = 20770.9995
$ws.Range("M126").Value = -11064.071
$ws.Range("N126").Value = -25710.9995
$ws.Range("H136").Value = 14372
$ws.Range("I136").Value = 3105.9
$ws.Range("K136").Value = 9317.700000000001
$ws.Range("M136").Value = -6767.700000000001
$ws.Range("H138").Value = 83180.73
$ws.Range("I138").Value = 52497
$ws.Range("J138").Value = 89999.336
$ws.Range("K138").Value = 52497
$ws.Range("L138").Value = 89999.336
$ws.Range("M138").Value = -47357
$ws.Range("N138").Value = -100279.336
